# Apply updated dSF (column F) values as part of a data re-pull / mean
# recalculation pass. Only column F values change; all other columns and
# rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -3
    6  = -5
    7  = -8
    11 = 2
    12 = 0
    13 = -5
    15 = -1
    20 = -4
    22 = -3
    24 = -3
    25 = -1
    26 = 1
    28 = -1
    31 = 4
    32 = 3
    34 = -1
    43 = 4
    46 = 3
    48 = -4
    49 = 3
    50 = -1
    58 = 0
    61 = 6
    62 = 0
    67 = 6
    68 = 11
    70 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
